$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "157×6=" "391×5="
Replace-Text "632×2=" "207×9="
Replace-Text "924×7=" "690×3="
Replace-Text "660×4=" "737×4="
Replace-Text "406×2=" "596×8="
Replace-Text "646×2=" "104×4="
Replace-Text "615×5=" "622×7="
Replace-Text "453×3=" "640×9="
Replace-Text "293×6=" "174×5="
Replace-Text "694×8=" "705×3="
Replace-Text "347×2=" "483×6="
Replace-Text "360×7=" "196×7="
Replace-Text "823×3=" "952×6="
Replace-Text "788×6=" "213×5="
Replace-Text "878×3=" "580×8="
Replace-Text "886×3=" "914×3="
Replace-Text "730×2=" "510×7="
Replace-Text "907×4=" "101×8="
Replace-Text "225×3=" "424×7="
Replace-Text "309×9=" "796×2="
Replace-Text "900×6=" "245×9="
Replace-Text "841×2=" "934×9="
Replace-Text "881×3=" "697×7="
Replace-Text "415×4=" "668×6="
Replace-Text "842×8=" "782×8="

Write-Host "Done"
